$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; existing rows 7-15 shift down to 8-16.
$ws.Rows("7:7").Insert()

# Fill the newly-inserted row 7 with the new weekly record.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44453
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112045
$ws.Range("G7").Value = "Zapallo"
$ws.Range("H7").Value = "Camote"
$ws.Range("I7").Value = "1a nueva(o)"
$ws.Range("J7").Value = 800
$ws.Range("K7").Value = 630
$ws.Range("L7").Value = 650
$ws.Range("M7").Value = 640
$ws.Range("N7").Value = "$/kilo (volumen en unidades)"
$ws.Range("O7").Value = "Perú"
$ws.Range("P7").Value = 640
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
